# Generate Report for handback
# - Overview/zh-cn/de-de status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# - zh-cn / de-de: populate "Latest Target File" (E) and "Latest Handback File" (F)
#   for rows 2 & 3, mirroring the Source File Name / Latest Handoff File hyperlinks.
# - zh-cn / de-de: stamp "Latest Handback DateTime" (G) for rows 2 & 3.

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet: refresh status column for both languages ---
$wsOverview.Range("B2").Value = $handedBack
$wsOverview.Range("C2").Value = $handedBack
$wsOverview.Range("B3").Value = $handedBack
$wsOverview.Range("C3").Value = $handedBack

# --- zh-cn sheet ---
$wsZhCn.Range("B2").Value = $handedBack
$wsZhCn.Range("B3").Value = $handedBack

# Row 2: 29a97ec9-0cb8-4386-afa4-0fc85831dfb1.md
$wsZhCn.Range("E2").Value = "29a97ec9-0cb8-4386-afa4-0fc85831dfb1.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/bb152afc2842070b38d3309a63050b4506b4653a/e2e/29a97ec9-0cb8-4386-afa4-0fc85831dfb1.md", "", "", "29a97ec9-0cb8-4386-afa4-0fc85831dfb1.md") | Out-Null

$wsZhCn.Range("F2").Value = "29a97ec9-0cb8-4386-afa4-0fc85831dfb1.75ce4d59c7705585df2d6be21ace3e95f5c70248.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4eaeea4fa144661f57987b23cf0b579d65120401/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/29a97ec9-0cb8-4386-afa4-0fc85831dfb1.75ce4d59c7705585df2d6be21ace3e95f5c70248.zh-cn.xlf", "", "", "29a97ec9-0cb8-4386-afa4-0fc85831dfb1.75ce4d59c7705585df2d6be21ace3e95f5c70248.zh-cn.xlf") | Out-Null

$wsZhCn.Range("G2").Value = "2016-02-16 10:36:57"

# Row 3: f12bf0be-48f2-4685-b862-21f679bd3c82.md
$wsZhCn.Range("E3").Value = "f12bf0be-48f2-4685-b862-21f679bd3c82.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/bb152afc2842070b38d3309a63050b4506b4653a/e2e/f12bf0be-48f2-4685-b862-21f679bd3c82.md", "", "", "f12bf0be-48f2-4685-b862-21f679bd3c82.md") | Out-Null

$wsZhCn.Range("F3").Value = "f12bf0be-48f2-4685-b862-21f679bd3c82.0a0e327d79a34b24fb16b690171ceeb7c34267b1.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4eaeea4fa144661f57987b23cf0b579d65120401/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f12bf0be-48f2-4685-b862-21f679bd3c82.0a0e327d79a34b24fb16b690171ceeb7c34267b1.zh-cn.xlf", "", "", "f12bf0be-48f2-4685-b862-21f679bd3c82.0a0e327d79a34b24fb16b690171ceeb7c34267b1.zh-cn.xlf") | Out-Null

$wsZhCn.Range("G3").Value = "2016-02-16 10:36:57"

# --- de-de sheet ---
$wsDeDe.Range("B2").Value = $handedBack
$wsDeDe.Range("B3").Value = $handedBack

# Row 2: 29a97ec9-0cb8-4386-afa4-0fc85831dfb1.md
$wsDeDe.Range("E2").Value = "29a97ec9-0cb8-4386-afa4-0fc85831dfb1.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/bb152afc2842070b38d3309a63050b4506b4653a/e2e/29a97ec9-0cb8-4386-afa4-0fc85831dfb1.md", "", "", "29a97ec9-0cb8-4386-afa4-0fc85831dfb1.md") | Out-Null

$wsDeDe.Range("F2").Value = "29a97ec9-0cb8-4386-afa4-0fc85831dfb1.75ce4d59c7705585df2d6be21ace3e95f5c70248.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20a7487ba3a8547531761fba62bfa504ff962688/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/29a97ec9-0cb8-4386-afa4-0fc85831dfb1.75ce4d59c7705585df2d6be21ace3e95f5c70248.de-de.xlf", "", "", "29a97ec9-0cb8-4386-afa4-0fc85831dfb1.75ce4d59c7705585df2d6be21ace3e95f5c70248.de-de.xlf") | Out-Null

$wsDeDe.Range("G2").Value = "2016-02-16 10:37:25"

# Row 3: f12bf0be-48f2-4685-b862-21f679bd3c82.md
$wsDeDe.Range("E3").Value = "f12bf0be-48f2-4685-b862-21f679bd3c82.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/bb152afc2842070b38d3309a63050b4506b4653a/e2e/f12bf0be-48f2-4685-b862-21f679bd3c82.md", "", "", "f12bf0be-48f2-4685-b862-21f679bd3c82.md") | Out-Null

$wsDeDe.Range("F3").Value = "f12bf0be-48f2-4685-b862-21f679bd3c82.0a0e327d79a34b24fb16b690171ceeb7c34267b1.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20a7487ba3a8547531761fba62bfa504ff962688/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f12bf0be-48f2-4685-b862-21f679bd3c82.0a0e327d79a34b24fb16b690171ceeb7c34267b1.de-de.xlf", "", "", "f12bf0be-48f2-4685-b862-21f679bd3c82.0a0e327d79a34b24fb16b690171ceeb7c34267b1.de-de.xlf") | Out-Null

$wsDeDe.Range("G3").Value = "2016-02-16 10:37:25"
